# Apply the cryptos.xlsx update: refreshed prices/volumes and a couple of
# row swaps (rank ties resolved differently this run).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = "'64.840.09"
$ws.Range('E2').Value = "'  -0.41%  "

# Row 3
$ws.Range('D3').Value = "'3.515.68"
$ws.Range('E3').Value = "'  -1.58%  "

# Row 4
$ws.Range('E4').Value = "'  +0.06%  "

# Row 5
$ws.Range('D5').Value = "'595.82"
$ws.Range('E5').Value = "'  -0.34%  "

# Row 6
$ws.Range('D6').Value = "'133.61"
$ws.Range('E6').Value = "'  -3.40%  "

# Row 7
$ws.Range('D7').Value = "'3.515.07"
$ws.Range('E7').Value = "'  -1.71%  "

# Row 9
$ws.Range('D9').Value = "'0.492"
$ws.Range('E9').Value = "'  -0.76%  "

# Row 10
$ws.Range('D10').Value = "'0.124"
$ws.Range('E10').Value = "'  -0.02%  "

# Row 11
$ws.Range('D11').Value = "'7.16"
$ws.Range('E11').Value = "'  +2.18%  "

# Row 12
$ws.Range('D12').Value = "'0.382"
$ws.Range('E12').Value = "'  -1.64%  "

# Row 13
$ws.Range('D13').Value = "'4.118.11"
$ws.Range('E13').Value = "'  -1.30%  "

# Row 14
$ws.Range('D14').Value = "'27.53"
$ws.Range('E14').Value = "'  +0.05%  "

# Row 15
$ws.Range('D15').Value = "'0.0000181"
$ws.Range('E15').Value = "'  -1.75%  "

# Row 16
$ws.Range('E16').Value = "'  +0.09%  "

# Row 17
$ws.Range('D17').Value = "'3.517.96"
$ws.Range('E17').Value = "'  -1.30%  "

# Row 18
$ws.Range('D18').Value = "'64.898.44"
$ws.Range('E18').Value = "'  -0.08%  "

# Row 19
$ws.Range('E19').Value = "'  -1.77%  "

# Row 20
$ws.Range('D20').Value = "'14.29"
$ws.Range('E20').Value = "'  -0.84%  "

# Row 21
$ws.Range('D21').Value = "'5.67"
$ws.Range('E21').Value = "'  -3.76%  "

# Row 22
$ws.Range('D22').Value = "'391.45"
$ws.Range('E22').Value = "'  -0.40%  "

# Row 23
$ws.Range('D23').Value = "'0.576"
$ws.Range('E23').Value = "'  -0.65%  "

# Row 24
$ws.Range('D24').Value = "'3.659.49"
$ws.Range('E24').Value = "'  -1.31%  "

# Row 25
$ws.Range('D25').Value = "'73.97"
$ws.Range('E25').Value = "'  -0.19%  "

# Row 26
$ws.Range('D26').Value = "'1.00"
$ws.Range('E26').Value = "'  +0.11%  "

# Row 27
$ws.Range('D27').Value = "'0.0000111"
$ws.Range('E27').Value = "'  -4.77%  "

# Row 28
$ws.Range('D28').Value = "'7.62"
$ws.Range('E28').Value = "'  -2.26%  "

# Row 29
$ws.Range('E29').Value = "'  +8.76%  "

# Row 30
$ws.Range('D30').Value = "'0.995"
$ws.Range('E30').Value = "'  -0.42%  "

# Row 31
$ws.Range('D31').Value = "'2.28"
$ws.Range('E31').Value = "'  -0.57%  "

# Row 32
$ws.Range('D32').Value = "'8.29"
$ws.Range('E32').Value = "'  -0.62%  "

# Row 33
$ws.Range('D33').Value = "'3.521.48"
$ws.Range('E33').Value = "'  -1.57%  "

# Row 34
$ws.Range('D34').Value = "'24.15"
$ws.Range('E34').Value = "'  +0.18%  "

# Row 35
$ws.Range('E35').Value = "'  +0.03%  "

# Row 36
$ws.Range('E36').Value = "'  -0.64%  "

# Row 37
$ws.Range('B37').Value = "NEARProtocol"
$ws.Range('C37').Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range('D37').Value = "'5.22"
$ws.Range('E37').Value = "'  +3.11%  "

# Row 38
$ws.Range('B38').Value = "ImmutableX"
$ws.Range('C38').Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range('D38').Value = "'1.57"
$ws.Range('E38').Value = "'  -0.02%  "

# Row 39
$ws.Range('B39').Value = "Monero"
$ws.Range('C39').Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range('D39').Value = "'168.48"
$ws.Range('E39').Value = "'  -0.56%  "

# Row 40
$ws.Range('D40').Value = "'6.82"
$ws.Range('E40').Value = "'  -1.61%  "

# Row 41
$ws.Range('D41').Value = "'0.0816"
$ws.Range('E41').Value = "'  +0.22%  "

# Row 42
$ws.Range('D42').Value = "'0.823"
$ws.Range('E42').Value = "'  -0.76%  "

# Row 43
$ws.Range('B43').Value = "EnergySwap"
$ws.Range('C43').Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range('D43').Value = "'25.81"
$ws.Range('E43').Value = "'  -3.02%  "

# Row 44
$ws.Range('B44').Value = "ONDO"
$ws.Range('C44').Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range('D44').Value = "'1.24"
$ws.Range('E44').Value = "'  +0.66%  "

# Row 45
$ws.Range('B45').Value = "OKB"
$ws.Range('C45').Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range('D45').Value = "'42.71"
$ws.Range('E45').Value = "'  -0.82%  "

# Row 46
$ws.Range('E46').Value = "'  -0.01%  "

# Row 47
$ws.Range('D47').Value = "'4.41"
$ws.Range('E47').Value = "'  -1.54%  "

# Row 48
$ws.Range('E48').Value = "'  -2.56%  "

# Row 49
$ws.Range('D49').Value = "'6.89"
$ws.Range('E49').Value = "'  -0.65%  "

# Row 50
$ws.Range('D50').Value = "'2.376.11"
$ws.Range('E50').Value = "'  -4.41%  "

# Row 51
$ws.Range('E51').Value = "'  +0.36%  "
